$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shape = $s.Shapes.Item(1)
$table = $shape.Table
Write-Host $table.Style
